$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.366.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.179.42"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.84"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.95%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.17%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.85"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.77"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.504.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.172.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.770"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.267.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.33%  "

# Row 19
$ws.Range("E19").Value = "  -0.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.57"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.61"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.29%  "

# Row 24
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.19%  "

# Row 25
$ws.Range("E25").Value = "  -0.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.40%  "

# Row 28
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +13.80%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.26"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0810"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.12"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.28%  "

# Row 35
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("E36").Value = "  +3.58%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0337"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.06"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.82%  "

# Row 41
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.196"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.81%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "59.60"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.37%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.16"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.62%  "

# Row 45
$ws.Range("E45").Value = "  +16.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.21%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.06%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0969"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.28%  "

# Row 49
$ws.Range("E49").Value = "  +0.59%  "

# Row 50
$ws.Range("E50").Value = "  +0.25%  "

# Row 51
$ws.Range("E51").Value = "  +0.63%  "
